# MasterFile.xlsx update:
#  - The worksheet link for the "Light" lesson (row 7) was filed under the
#    wrong column (D, "Quiz") and is moved to the correct column (E, "Worksheet").
#  - The worksheet link for the "Forest Our Life Line" lesson (row 8, column E)
#    is updated to a more specific path.
#  - The active selection is moved to E12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "/Light.pdf" from D7 (wrong column) to E7 (Worksheet column).
$ws.Range("E7").Value = $ws.Range("D7").Value2
$ws.Range("D7").Value = ""

# Update the forest worksheet link to the full path.
$ws.Range("E8").Value = "/pdfs/grade7/science/forest.pdf"

# Update the active cell / selection shown in the saved view.
$ws.Range("E12").Select()
